$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write literal text into a cell while keeping it a genuine shared
# string (t="s") instead of letting Excel's type-inference turn things like
# "01/01/2012" into a real date serial, and without leaving extra unused
# style entries behind. Pattern: paste the formatting from a reference cell
# that already carries the desired style, push the text in through a
# quoted formula (forces text), then immediately convert that formula back
# down to a plain value with Paste Values.
function Set-TextCell($addr, $text, $styleRefAddr) {
    $ws.Range($styleRefAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
    $ws.Range($addr).Copy()
    $excel.CutCopyMode = $false
}

# Row 10: "Objetivos:" value replaced by the teacher info text
Set-TextCell "B10" "7455355 - Robson da Silva Rocha" "B10"
Set-TextCell "C10" "7455355 - Robson da Silva Rocha" "C10"

# Row 13: previously a blank-label row holding the teacher info; now becomes
# "Programa resumido:" with value "Semestral"
Set-TextCell "A13" "Programa resumido:" "A3"
Set-TextCell "B13" "Semestral" "B9"
Set-TextCell "C13" "Semestral" "C9"

# Row 14: label shifts to "Short syllabus:"; its long text value is removed
Set-TextCell "A14" "Short syllabus:" "A3"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# Row 15: label shifts to "Programa:"; gains the activation date value
Set-TextCell "A15" "Programa:" "A3"
Set-TextCell "B15" "01/01/2012" "B9"
Set-TextCell "C15" "01/01/2012" "C9"

# Row 16: label shifts to "Syllabus:"; its long text value is removed
Set-TextCell "A16" "Syllabus:" "A3"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# Row 17: label shifts to "Avaliação:"
Set-TextCell "A17" "Avaliação:" "A3"

# Row 18: label shifts to "Método:"; gains the teacher info value
Set-TextCell "A18" "Método:" "A3"
Set-TextCell "B18" "7455355 - Robson da Silva Rocha" "B9"
Set-TextCell "C18" "7455355 - Robson da Silva Rocha" "C9"

# Row 19: label shifts to "Critério:" (B/C keep their existing text)
Set-TextCell "A19" "Critério:" "A3"

# Row 20: label shifts to "Norma de recuperação:" (B/C keep their existing text)
Set-TextCell "A20" "Norma de recuperação:" "A3"

# Row 21: label shifts to "Bibliografia:" (B/C keep their existing text)
Set-TextCell "A21" "Bibliografia:" "A3"

# Row 22 (old Bibliografia long text row) is removed entirely
$ws.Rows(22).Delete()

# Fix up row heights to match the new layout
$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
